# Registers.xlsx - Add files via upload
# DecisionRegister: decision #2 ("Change of project scope") gets approved
# (Approved By / Date filled in) and a new decision #3 ("Use OptiTrack to
# sense the match") is added to the Decision table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DecisionRegister")

# --- Fill in the "Approved By" / "Date" columns for the existing row (ID 2) ---
$ws.Range("G3").Value = "Team (majority vote), René and Erjen"
$ws.Range("H3").Value = 45353   # 3/2/2024

# --- Append a new row to the Decision table (Table2) ---
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Copy the number/status formatting from the row above so the new cells
# keep the same date format / "Good" style used throughout the table.
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H4").PasteSpecial(-4122)

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Use OptiTrack to sense the match"
$ws.Range("C4").Value = "What data the solutions uses"
$ws.Range("D4").Value = 45355   # 3/4/2024
$ws.Range("E4").Value = "Team"
$ws.Range("F4").Value = "Not approved"

# --- Make DecisionRegister the active sheet / selected cell ---
$ws.Activate()
$ws.Range("F5").Select()
